$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.366.63"
$ws.Range("D3").Value = "1.592.12"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "1.816.46"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "1.594.54"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "26.381.14"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "212.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.63%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").Value = "1.337.18"
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -20.13%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.767"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "1.729.22"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0986"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.47%  "
